{"js": "// Add two new checklist bullet items after the \"PHP.ini Go to line no 916 ...\"\n// item: \"PHP.ini Enable open_ssl port\" and\n// \"No Echo Statement in any of the Controller or Model\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the last paragraph that contains the anchor text so the new\n// bullets are inserted right after it, inheriting its list formatting.\nconst items = paragraphs.items;\nlet anchor = items[items.length - 1];\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"PHP.ini Go to line no 916\") !== -1) {\n    anchor = items[i];\n    break;\n  }\n}\n\nconst firstNew = anchor.insertParagraph(\n  \"PHP.ini Enable open_ssl port\",\n  Word.InsertLocation.after\n);\nfirstNew.insertParagraph(\n  \"No Echo Statement in any of the Controller or Model\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Add two new checklist bullet items after the \"PHP.ini Go to line no 916 ...\"\n# item: \"PHP.ini Enable open_ssl port\" and\n# \"No Echo Statement in any of the Controller or Model\".\n$d = $word.ActiveDocument\n\n# Find the index of the anchor paragraph containing the PHP.ini\n# upload_max_filesize text so the new bullets land right after it,\n# inheriting its list formatting.\n$anchorIndex = $d.Paragraphs.Count\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*PHP.ini Go to line no 916*\") {\n        $anchorIndex = $i\n    }\n}\n\n$anchor = $d.Paragraphs.Item($anchorIndex)\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Text = \"PHP.ini Enable open_ssl port\"\n\n$second = $d.Paragraphs.Item($anchorIndex + 1)\n$second.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($anchorIndex + 2).Range.Text = \"No Echo Statement in any of the Controller or Model\"\n"}
